$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18 ---
# A18, B18, D18 unchanged (NAME, CONDITION, ACTION)
$ws.Cells.Item(18,3).Value2 = "ACTION"
$ws.Cells.Item(18,5).Value2 = "ACTION"

# --- Row 19 ---
$ws.Cells.Item(19,2).Value2 = "Senior perk"
$ws.Cells.Item(19,3).Value2 = "Test"
$ws.Cells.Item(19,4).Value2 = "jnn"
$ws.Cells.Item(19,5).Value2 = "jn"

# --- Row 20 ---
$ws.Cells.Item(20,3).Value2 = "gfndnvbx"
$ws.Cells.Item(20,4).ClearContents()
$ws.Cells.Item(20,5).ClearContents()

# --- Row 21 ---
$ws.Cells.Item(21,3).Value2 = "dsgagass"
$ws.Cells.Item(21,4).ClearContents()
$ws.Cells.Item(21,5).ClearContents()

# --- Row 22 ---
$ws.Cells.Item(22,3).Value2 = "dsvsbsb"
$ws.Cells.Item(22,4).ClearContents()
$ws.Cells.Item(22,5).ClearContents()

# --- Row 23 ---
$ws.Cells.Item(23,3).Value2 = "nsngnsg"
$ws.Cells.Item(23,4).ClearContents()
$ws.Cells.Item(23,5).ClearContents()

# --- Row 24 ---
$ws.Cells.Item(24,3).Value2 = "Code changed 10010018"
$ws.Cells.Item(24,4).ClearContents()
$ws.Cells.Item(24,5).ClearContents()

# --- Row 25 ---
$ws.Cells.Item(25,1).Value2 = "Test"
$ws.Cells.Item(25,2).ClearContents()
$ws.Cells.Item(25,3).ClearContents()
$ws.Cells.Item(25,4).ClearContents()
$ws.Cells.Item(25,5).ClearContents()

# --- Row 26 ---
$ws.Cells.Item(26,1).Value2 = "Test"
$ws.Cells.Item(26,2).ClearContents()
$ws.Cells.Item(26,3).ClearContents()
$ws.Cells.Item(26,4).ClearContents()
$ws.Cells.Item(26,5).ClearContents()

# --- Row 27 removed entirely ---
$ws.Rows("27:27").Delete()
